# Actualización automática del tracker
# Fill in the pending result/profit for matches that have since concluded,
# and append the newest fixture row to the bottom of the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resolve pending rows (previously blank "resultado"/"profit" columns)
$ws.Range("G114").Value = "Fallo"
$ws.Range("H114").Value = -1

$ws.Range("G116").Value = "Acierto"
$ws.Range("H116").Value = 0.62

$ws.Range("G117").Value = "Fallo"
$ws.Range("H117").Value = -1

$ws.Range("G118").Value = "Fallo"
$ws.Range("H118").Value = -1

$ws.Range("G119").Value = "Fallo"
$ws.Range("H119").Value = -1

# Append the new fixture as row 121 (still pending -> resultado/profit blank)
$ws.Range("A121").Value = 14807075

# Force the date to be stored as literal text (e.g. "2025-10-08") instead of
# letting Excel auto-convert the recognizable date string into a date serial.
$ws.Range("B121").NumberFormat = "@"
$ws.Range("B121").Value = "2025-10-08"
$ws.Range("B121").ClearFormats()

$ws.Range("C121").Value = "Niels Visker"
$ws.Range("D121").Value = "Ryan Peniston"
$ws.Range("E121").Value = "Gana Niels Visker"
$ws.Range("F121").Value = 3.75

# resultado / profit are still pending for this new fixture -- keep the
# cells present (matching the rest of the tracker's row shape) but empty.
$ws.Range("G121").Style = "Normal"
$ws.Range("H121").Style = "Normal"
